$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.952.42"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "1.892.96"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'0.8207"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +6.88%  "
$ws.Range("D6").Value = "'241.33"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.3220"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +5.60%  "
$ws.Range("D9").Value = "'26.46"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +4.40%  "
$ws.Range("D10").Value = "'0.07012"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.48%  "
$ws.Range("D11").Value = "'0.08033"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Value = "'0.7467"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("D13").Value = "1.899.36"
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").Value = "'5.192"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "'92.20"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").Value = "29.944.62"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "'14.02"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").Value = "'5.888"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "'244.77"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "'0.000007744"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "2.138.49"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").Value = "'1.004"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").Value = "'6.896"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").Value = "'0.1594"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +24.13%  "
$ws.Range("D26").Value = "'166.46"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").Value = "'9.177"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("D28").Value = "'18.84"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").Value = "'2.068"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.06%  "
$ws.Range("D30").Value = "'1.366"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("D31").Value = "'1.516"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("D32").Value = "'4.261"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("D33").Value = "'0.05622"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +7.19%  "
$ws.Range("D34").Value = "'4.069"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").Value = "'1.272"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.51%  "
$ws.Range("D36").Value = "'0.7303"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("D37").Value = "'2.723"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").Value = "'0.01913"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").Value = "'2.783"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").Value = "'0.4406"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "'71.86"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").Value = "'5.945"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.02%  "
$ws.Range("D43").Value = "'0.8429"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.00%  "
$ws.Range("D45").Value = "'1.879"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").Value = "'7.570"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "'100.58"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.679"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'989.51"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +9.10%  "
$ws.Range("D50").Value = "2.041.13"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("D51").Value = "'35.95"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.61%  "